# Update the User Master workbook:
#  - Change the "Password" column (C) values for all data rows from the
#    numeric placeholder 123 to the text value "abc123" so the repository's
#    updated password-comparison logic can be exercised against a realistic,
#    non-numeric password string.
#  - Leave the cursor/selection on D5, matching where the user ended up
#    after editing the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User")

# Data rows are 2-11, Password is column C.
$ws.Range("C2:C11").Value = "abc123"

# Update the active selection to D5 (as left by the editor).
$ws.Range("D5").Select() | Out-Null
